# Nakka_LabExam03Grading.xlsx - "kalyankar to pusapati done"
# Fill in the "Points for grading" (column E) scores for the Generic
# section (rows 3-6) and the Customer Class section (rows 10-14) to
# mirror the "Total Points" already recorded in column D, then leave
# the selection on E15 (the Customer Class section total).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

$ws.Range("E15").Select()
